# Regenerate the K (column G) values for each row, replacing the old
# Strike#-derived values with the new computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 1
    4  = 0
    5  = 2
    6  = 1
    7  = 1
    8  = 2
    9  = 3
    10 = 0
    11 = 7
    12 = 3
    13 = 2
    14 = 5
    15 = 2
    16 = 5
    17 = 3
    18 = 5
    19 = 4
    20 = 5
    21 = 6
    22 = 0
    23 = 3
    24 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
